$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in / correct the dates for the Feature 2 / Feature 3 / Feature 4 rows ---
# These rows previously had placeholder dates (18/05/2024, 18/05/2025) that are
# corrected to the actual date (18/05/2023), and the "Feature 4: collapsable
# controls" row (15) gets its date filled in for the first time.
$ws.Range("C10").Value = "18/05/2023"
$ws.Range("C11").Value = "18/05/2023"
$ws.Range("C12").Value = "18/05/2023"
$ws.Range("C13").Value = "18/05/2023"
$ws.Range("C14").Value = "18/05/2023"
$ws.Range("C15").Value = "18/05/2023"

# Row 15 grew taller once the date + wrapped content settled.
$ws.Rows("15:15").RowHeight = 109.2

# --- Insert a brand-new log entry (row 16) for the "Animate the algorithm" task ---
$ws.Rows("16:16").Insert()

# Copy the formatting from the blank spacer row (now row 17, right below the
# freshly inserted row) onto the new row so it matches the rest of the table.
$ws.Range("A17:F17").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)

$ws.Range("A16").Value = "Feature Extra: Animate the algorithm"
$ws.Range("B16").Value = 1.25
$ws.Range("C16").Value = "18/05/2023"
$ws.Range("D16").Value = "I have added an additional task to implement animated maze generation.`nTo accomplish this, I have created a new function that returns animation frames, represented as a boolean array per frame. Storing all the frames should not be a concern since the array is a byte array, which doesn't consume excessive memory.`nHowever, I noticed a previous issue where all prefabs were being deleted when creating a new maze. To address this problem and avoid excessive memory allocation and release, I have made some changes. Now, I create an object pool of prefabs based on the size of the first frame. These prefabs are initially disabled in the scene. During the animation, I simply enable or disable the prefabs based on whether they should represent a wall or a path in the frame."
$ws.Range("E16").Value = "X"

$ws.Rows("16:16").RowHeight = 157.2
$ws.Rows("17:17").RowHeight = 16.8

# --- Update the view so the newly added row is visible / selected ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B17").Select()
